$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-10: I = 1, J = same value as H in that row
$hValues = @{
    2 = 7
    3 = 6
    4 = 8
    5 = 6
    6 = 2
    7 = 4
    8 = 4
    9 = 6
    10 = 7
}

foreach ($r in $hValues.Keys) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hValues[$r]
}

# Data rows 11-12: I = 8, J = 8
foreach ($r in 11..12) {
    $ws.Cells.Item($r, 9).Value = 8
    $ws.Cells.Item($r, 10).Value = 8
}
